$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Trial 37", 0, 0, 0, 80.53076314926147, 81.70234942436218, 86.79865407943726, 81.36654353141785, 36.32433462142944, 3.009976148605347),
    @("Trial 38", 0, 0, 0, 663.6075005531311, 666.0170459747314, 683.7945878505707, 664.8327276706696, 289.1930079460144, 14.67776441574097)
)

$startRow = 38
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
